# Purchase Request workbook update — March 27, 2020 commit
# "Links, Fonts, redirect"
#
# Re-targets the PR from the Office of the Regional Director (ORD) to
# LGMED for a "Meeting of the RTF ELCAC": swaps the single line item from
# the Regional Team Conference venue/meals package to a Meals-only line,
# updates the requested-by officer's printed name/designation, and
# updates the quantity / unit cost / total of that line item. The Grand
# Total cell (F36) keeps its existing formula =SUM(F11:F35) and is left
# for Excel to recompute from the new totals.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header block: Office/Section + PR No. -----------------------------
$ws.Range("B7").Value = "LGMED"
$ws.Range("C7").Value = "PR No.:  2020-03-0160"

# --- Line item (row 11) --------------------------------------------------
$ws.Range("A11").Value = "S863"
$ws.Range("B11").Value = "pax"
$ws.Range("C11").Value = "Meals (AM Snack-Lunch-PM Snacks)`n"
$ws.Range("D11").Value = 110
$ws.Range("E11").Value = 500
$ws.Range("F11").Value = 55000

# --- Purpose --------------------------------------------------------------
$ws.Range("B37").Value = "Meeting of the RTF ELCAC"

# --- Requested-by signatory ----------------------------------------------
$ws.Range("B43").Value = "GILBERTO L. TUMAMAC"
$ws.Range("B44").Value = "OIC - LGMED Chief"
